$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.044.44'
$ws.Range('E2').Value = '  +5.13%  '
$ws.Range('D3').Value = '2.285.94'
$ws.Range('E3').Value = '  +2.09%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''231.05'
$ws.Range('E5').Value = '  -0.80%  '
$ws.Range('D6').Value = '''0.624'
$ws.Range('E6').Value = '  -0.71%  '
$ws.Range('E7').Value = '  -1.18%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '''0.425'
$ws.Range('E9').Value = '  +4.31%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').Value = '''57.90'
$ws.Range('E10').Value = '  -1.88%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = '''0.0940'
$ws.Range('E11').Value = '  +1.95%  '
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('D13').Value = '2.630.79'
$ws.Range('E13').Value = '  +2.57%  '
$ws.Range('D14').Value = '''24.26'
$ws.Range('E14').Value = '  +8.33%  '
$ws.Range('D15').Value = '''15.71'
$ws.Range('E15').Value = '  -0.32%  '
$ws.Range('D16').Value = '''5.90'
$ws.Range('E16').Value = '  +4.70%  '
$ws.Range('D17').Value = '''0.813'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('D18').Value = '2.302.38'
$ws.Range('E18').Value = '  +2.41%  '
$ws.Range('D19').Value = '43.917.21'
$ws.Range('E19').Value = '  +5.11%  '
$ws.Range('D20').Value = '0.0₃0943'
$ws.Range('E20').Value = '  +2.80%  '
$ws.Range('D21').Value = '''73.70'
$ws.Range('E21').Value = '  +1.37%  '
$ws.Range('D22').Value = '''6.25'
$ws.Range('E22').Value = '  +2.76%  '
$ws.Range('D23').Value = '''250.72'
$ws.Range('E23').Value = '  -1.58%  '
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('E25').Value = '  +6.43%  '
$ws.Range('E26').Value = '  +2.34%  '
$ws.Range('D27').Value = '''9.88'
$ws.Range('E27').Value = '  +1.89%  '
$ws.Range('D28').Value = '''171.15'
$ws.Range('E28').Value = '  +0.91%  '
$ws.Range('D29').Value = '''0.139'
$ws.Range('E29').Value = '  -4.03%  '
$ws.Range('D30').Value = '''20.55'
$ws.Range('E30').Value = '  +1.96%  '
$ws.Range('D31').Value = '''1.42'
$ws.Range('E31').Value = '  -0.51%  '
$ws.Range('D32').Value = '''2.75'
$ws.Range('E32').Value = '  +1.06%  '
$ws.Range('E33').Value = '  -0.49%  '
$ws.Range('D34').Value = '''4.80'
$ws.Range('E34').Value = '  +2.07%  '
$ws.Range('D35').Value = '''5.05'
$ws.Range('E35').Value = '  -0.49%  '
$ws.Range('D36').Value = '''0.0658'
$ws.Range('E36').Value = '  +2.47%  '
$ws.Range('D37').Value = '''6.52'
$ws.Range('E37').Value = '  -2.29%  '
$ws.Range('D38').Value = '''3.64'
$ws.Range('E38').Value = '  -2.82%  '
$ws.Range('D39').Value = '''2.40'
$ws.Range('E39').Value = '  +1.24%  '
$ws.Range('D40').Value = '''0.0250'
$ws.Range('E40').Value = '  +3.47%  '
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('D42').Value = '''8.75'
$ws.Range('E42').Value = '  +0.78%  '
$ws.Range('D43').Value = '''0.000222'
$ws.Range('E43').Value = '  -14.88%  '
$ws.Range('D44').Value = '''0.0970'
$ws.Range('E44').Value = '  +0.98%  '
$ws.Range('B45').Value = 'Celestia'
$ws.Range('C45').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D45').Value = '''10.42'
$ws.Range('E45').Value = '  +17.85%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '''98.71'
$ws.Range('E46').Value = '  -1.19%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').Value = '''1.21'
$ws.Range('E47').Value = '  -1.75%  '
$ws.Range('D48').Value = '''17.04'
$ws.Range('E48').Value = '  +2.52%  '
$ws.Range('B49').Value = 'FTXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D49').Value = '''4.40'
$ws.Range('E49').Value = '  -6.01%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '1.478.30'
$ws.Range('E50').Value = '  -0.58%  '
$ws.Range('D51').Value = '''1.09'
$ws.Range('E51').Value = '  +0.29%  '
